$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 550
$ws.Range("I7").Value = 550
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 550
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -438
$ws.Range("N7").ClearContents()
$ws.Range("H14").Value = 550
$ws.Range("I14").Value = 550
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 550
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -359
$ws.Range("N14").ClearContents()
$ws.Range("H43").Value = 4999.5
$ws.Range("I43").Value = 6999
$ws.Range("J43").Value = 3000
$ws.Range("K43").Value = 6999
$ws.Range("L43").Value = 3000
$ws.Range("M43").Value = -6930
$ws.Range("N43").Value = -3138
$ws.Range("H47").Value = 16344.333
$ws.Range("I47").Value = 17813.2
$ws.Range("K47").Value = 17813.2
$ws.Range("M47").Value = -16841.2
$ws.Range("H134").Value = 95390
$ws.Range("J134").Value = 95390
$ws.Range("L134").Value = 95390
$ws.Range("N134").Value = -105530

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 881.25
$ws.Range("I110").Value = 904
$ws.Range("K110").Value = 904
$ws.Range("M110").Value = 1141

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1486.0769
$ws.Range("I20").Value = 1482.8572
$ws.Range("K20").Value = 1482.8572
$ws.Range("M20").Value = -1235.8572
$ws.Range("H99").Value = 2500
$ws.Range("I99").Value = 2500
$ws.Range("K99").Value = 2500
$ws.Range("M99").Value = -1002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 93750220
$ws.Range("I7").Value = 45454820
$ws.Range("K7").Value = 45454820
$ws.Range("M7").Value = -45454707
$ws.Range("H12").Value = 388
$ws.Range("I12").Value = 388
$ws.Range("K12").Value = 388
$ws.Range("M12").Value = -218
$ws.Range("H62").Value = 35722640
$ws.Range("I62").Value = 41674748
$ws.Range("J62").Value = 9999
$ws.Range("K62").Value = 41674748
$ws.Range("L62").Value = 9999
$ws.Range("M62").Value = -41674124
$ws.Range("N62").Value = -11247
$ws.Range("H65").Value = 35722640
$ws.Range("I65").Value = 41674748
$ws.Range("J65").Value = 9999
$ws.Range("K65").Value = 208373740
$ws.Range("L65").Value = 49995
$ws.Range("M65").Value = -208370620
$ws.Range("N65").Value = -56235
$ws.Range("H69").Value = 8750
$ws.Range("I69").Value = 8750
$ws.Range("K69").Value = 8750
$ws.Range("M69").Value = -8001
$ws.Range("H72").Value = 8750
$ws.Range("I72").Value = 8750
$ws.Range("K72").Value = 26250
$ws.Range("M72").Value = -22506
$ws.Range("H93").Value = 14884.5
$ws.Range("I93").Value = 14884.5
$ws.Range("K93").Value = 14884.5
$ws.Range("M93").Value = -13012.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 10009.5
$ws.Range("I10").Value = 19
$ws.Range("K10").Value = 57
$ws.Range("M10").Value = 82
$ws.Range("H39").Value = 2395.4
$ws.Range("J39").Value = 5000
$ws.Range("L39").Value = 15000
$ws.Range("N39").Value = -15588
$ws.Range("H59").Value = 5000
$ws.Range("I59").Value = 5000
$ws.Range("K59").Value = 15000
$ws.Range("M59").Value = -14460
$ws.Range("H86").Value = 8386
$ws.Range("I86").Value = 675
$ws.Range("J86").Value = 18667.334
$ws.Range("K86").Value = 2025
$ws.Range("L86").Value = 56002.00199999999
$ws.Range("M86").Value = -839
$ws.Range("N86").Value = -58374.00199999999
$ws.Range("H89").Value = 8386
$ws.Range("I89").Value = 675
$ws.Range("J89").Value = 18667.334
$ws.Range("K89").Value = 6075
$ws.Range("L89").Value = 168006.006
$ws.Range("M89").Value = -147
$ws.Range("N89").Value = -179862.006
$ws.Range("H98").Value = 817.625
$ws.Range("I98").Value = 1115.25
$ws.Range("J98").Value = 520
$ws.Range("K98").Value = 3345.75
$ws.Range("L98").Value = 1560
$ws.Range("M98").Value = -1847.75
$ws.Range("N98").Value = -4556
$ws.Range("H129").Value = 3799.9
$ws.Range("J129").Value = 4166.5
$ws.Range("L129").Value = 12499.5
$ws.Range("N129").Value = -22499.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 47.46154
$ws.Range("I2").Value = 51.42857
$ws.Range("J2").Value = 42.833332
$ws.Range("K2").Value = 51.42857
$ws.Range("L2").Value = 42.833332
$ws.Range("M2").Value = 61.57143
$ws.Range("N2").Value = -268.833332
$ws.Range("H59").Value = 1250
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 1250
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 1250
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -2416
$ws.Range("H70").Value = 99999
$ws.Range("I70").Value = 99999
$ws.Range("K70").Value = 99999
$ws.Range("M70").Value = -99729
$ws.Range("H73").Value = 99999
$ws.Range("I73").Value = 99999
$ws.Range("K73").Value = 99999
$ws.Range("M73").Value = -99063
$ws.Range("H97").Value = 1234
$ws.Range("I97").Value = 1234
$ws.Range("K97").Value = 1234
$ws.Range("M97").Value = -738

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3306.0588
$ws.Range("I22").Value = 3071.2856
$ws.Range("J22").Value = 3470.4
$ws.Range("K22").Value = 3071.2856
$ws.Range("L22").Value = 3470.4
$ws.Range("M22").Value = -2776.2856
$ws.Range("N22").Value = -4060.4
$ws.Range("H27").Value = 3306.0588
$ws.Range("I27").Value = 3071.2856
$ws.Range("J27").Value = 3470.4
$ws.Range("K27").Value = 3071.2856
$ws.Range("L27").Value = 3470.4
$ws.Range("M27").Value = -2964.2856
$ws.Range("N27").Value = -3684.4
$ws.Range("H46").Value = 5002
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 5002
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 5002
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -5378
$ws.Range("H55").Value = 2874.6428
$ws.Range("I55").Value = 2468.5
$ws.Range("J55").Value = 3416.1667
$ws.Range("K55").Value = 2468.5
$ws.Range("L55").Value = 3416.1667
$ws.Range("M55").Value = -2295.5
$ws.Range("N55").Value = -3762.1667
$ws.Range("H68").Value = 6422.222
$ws.Range("I68").Value = 2714.2856
$ws.Range("J68").Value = 19400
$ws.Range("K68").Value = 2714.2856
$ws.Range("L68").Value = 19400
$ws.Range("M68").Value = -1965.2856
$ws.Range("N68").Value = -20898
$ws.Range("H71").Value = 6422.222
$ws.Range("I71").Value = 2714.2856
$ws.Range("J71").Value = 19400
$ws.Range("K71").Value = 13571.428
$ws.Range("L71").Value = 97000
$ws.Range("M71").Value = -9827.428
$ws.Range("N71").Value = -104488
$ws.Range("H100").Value = 5312.75
$ws.Range("J100").Value = 6300
$ws.Range("L100").Value = 6300
$ws.Range("N100").Value = -7382

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 30000
$ws.Range("I54").Value = 30000
$ws.Range("K54").Value = 30000
$ws.Range("M54").Value = -29480
$ws.Range("H75").Value = 75000
$ws.Range("J75").Value = 75000
$ws.Range("L75").Value = 75000
$ws.Range("N75").Value = -76872
$ws.Range("H78").Value = 75000
$ws.Range("J78").Value = 75000
$ws.Range("L78").Value = 225000
$ws.Range("N78").Value = -234360
$ws.Range("H81").Value = 2000.6666
$ws.Range("I81").Value = 2000.6666
$ws.Range("K81").Value = 4001.3332
$ws.Range("M81").Value = -2940.3332
$ws.Range("H84").Value = 2000.6666
$ws.Range("I84").Value = 2000.6666
$ws.Range("K84").Value = 20006.666
$ws.Range("M84").Value = -14702.666
